$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.105608157319738
$ws.Range("C2").Value = 0.235979606700655
$ws.Range("D2").Value = 0.0852148579752367
$ws.Range("E2").Value = 0.771303714493809
$ws.Range("F2").Value = 0.159504734158776
$ws.Range("G2").Value = 0.651857246904588
$ws.Range("H2").Value = 0.300072833211945
$ws.Range("I2").Value = 0.775673707210488
$ws.Range("J2").Value = 0.0305899490167516
$ws.Range("K2").Value = 0.777130371449381
$ws.Range("L2").Value = 0.91332847778587
$ws.Range("M2").Value = 0.0160233066278223
$ws.Range("N2").Value = 0.0320466132556446
$ws.Range("O2").Value = 0.353241077931537
$ws.Range("P2").Value = 0.0101966496722505
$ws.Range("Q2").Value = 0.0517115804806992
$ws.Range("R2").Value = 0.937363437727604
$ws.Range("S2").Value = 0.0546249089584851
$ws.Range("T2").Value = 0.599417334304443
$ws.Range("U2").Value = 0.191551347414421
$ws.Range("V2").Value = 0.0786598689002185
$ws.Range("W2").Value = 0.0830298616168973
$ws.Range("X2").Value = 0.117989803350328
$ws.Range("B3").Value = 0.836125273124545
$ws.Range("C3").Value = 0.0662782228696285
$ws.Range("D3").Value = 0.0466132556445739
$ws.Range("E3").Value = 0.155134741442098
$ws.Range("F3").Value = 0.0633648943918427
$ws.Range("G3").Value = 0.0436999271667881
$ws.Range("H3").Value = 0.0968681718863802
$ws.Range("I3").Value = 0.124544792425346
$ws.Range("J3").Value = 0.0356882738528769
$ws.Range("K3").Value = 0.0422432629278951
$ws.Range("L3").Value = 0.0247632920611799
$ws.Range("M3").Value = 0.0538965768390386
$ws.Range("N3").Value = 0.0487982520029133
$ws.Range("O3").Value = 0.00509832483612527
$ws.Range("P3").Value = 0.163146394756009
$ws.Range("Q3").Value = 0.00582665695557174
$ws.Range("R3").Value = 0.00145666423889294
$ws.Range("S3").Value = 0.225054624908958
$ws.Range("T3").Value = 0.0269482884195193
$ws.Range("U3").Value = 0.0400582665695557
$ws.Range("V3").Value = 0.0852148579752367
$ws.Range("W3").Value = 0.0772032046613256
$ws.Range("X3").Value = 0.0728332119446468
$ws.Range("B4").Value = 0.00728332119446468
$ws.Range("C4").Value = 0.412964311726147
$ws.Range("D4").Value = 0.823743627093955
$ws.Range("E4").Value = 0.0626365622723962
$ws.Range("F4").Value = 0.763292061179898
$ws.Range("G4").Value = 0.175528040786599
$ws.Range("H4").Value = 0.554989075018208
$ws.Range("I4").Value = 0.0750182083029862
$ws.Range("J4").Value = 0.92571012381646
$ws.Range("K4").Value = 0.106336489439184
$ws.Range("L4").Value = 0.0393299344501093
$ws.Range("M4").Value = 0.909686817188638
$ws.Range("N4").Value = 0.0626365622723962
$ws.Range("O4").Value = 0.640932265112891
$ws.Range("P4").Value = 0.0597232337946103
$ws.Range("Q4").Value = 0.000728332119446468
$ws.Range("R4").Value = 0.032774945375091
$ws.Range("S4").Value = 0.718135469774217
$ws.Range("T4").Value = 0.0138383102694829
$ws.Range("U4").Value = 0.753823743627094
$ws.Range("V4").Value = 0.0990531682447196
$ws.Range("W4").Value = 0.785870356882739
$ws.Range("X4").Value = 0.767662053896577
$ws.Range("B5").Value = 0.0509832483612527
$ws.Range("C5").Value = 0.284049526584122
$ws.Range("D5").Value = 0.0444282592862345
$ws.Range("E5").Value = 0.010924981791697
$ws.Range("F5").Value = 0.0131099781500364
$ws.Range("G5").Value = 0.128914785142025
$ws.Range("H5").Value = 0.0473415877640204
$ws.Range("I5").Value = 0.0247632920611799
$ws.Range("J5").Value = 0.00801165331391114
$ws.Range("K5").Value = 0.0742898761835397
$ws.Range("L5").Value = 0.0225782957028405
$ws.Range("M5").Value = 0.0196649672250546
$ws.Range("N5").Value = 0.856518572469046
$ws.Range("P5").Value = 0.766205389657684
$ws.Range("Q5").Value = 0.941005098324836
$ws.Range("R5").Value = 0.0276766205389658
$ws.Range("S5").Value = 0.0021849963583394
$ws.Range("T5").Value = 0.359796067006555
$ws.Range("U5").Value = 0.0145666423889294
$ws.Range("V5").Value = 0.737072104879825
$ws.Range("W5").Value = 0.0531682447195921
$ws.Range("X5").Value = 0.0407865986890022
